$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the FFR column (C) - this shifts the LF column (D) into C
$ws.Columns.Item(3).Delete()

# Remove the FFR Lag row (now row 3 after the column delete) - shifts the LF Lag row up
$ws.Rows.Item(3).Delete()

# Update the remaining coefficient values with the new Crisis and Credit Allocation figures.
# B2 ("0.27") looks numeric, so stage it as text in a scratch cell and paste-special the
# value in, keeping it a plain text cell (matching the source workbook) without altering
# B2's own number format/style.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "0.27"
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4163)  # xlPasteValues
$ws.Columns.Item(26).Delete()        # clean up the scratch column (Z)

$ws.Range("C2").Value = "-8.77*"
$ws.Range("B3").Value = "-0.11*"
$ws.Range("C3").Value = "2.45***"
